$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - SamsungElec (005930.KS) refreshed metrics
$ws.Range("D2").Value = 101100
$ws.Range("E2").Value = 46.7
$ws.Range("F2").Value = 4.55
$ws.Range("G2").Value = 50
$ws.Range("H2").Value = 70
$ws.Range("I2").Value = 70
$ws.Range("J2").Value = 83
$ws.Range("K2").Value = 68.8
$ws.Range("N2").Value = 85.92500513438651

# Row 3 - now 240810.KS ticker group
$ws.Range("B3").Value = "240810.KS,0P00017YB3,330568"
$ws.Range("C3").Value = "240810.KS"
$ws.Range("D3").Value = 61300
$ws.Range("E3").Value = 30.8
$ws.Range("F3").Value = 8.109999999999999
$ws.Range("G3").Value = 20
$ws.Range("H3").Value = 63
$ws.Range("I3").Value = 70
$ws.Range("J3").Value = 66
$ws.Range("K3").Value = 59.8
$ws.Range("M3").Value = "⛔ 관망하십시오."
$ws.Range("N3").Value = 85.92500513438651

# Row 4 - now DB HiTek (000990.KS)
$ws.Range("B4").Value = "DB HiTek"
$ws.Range("C4").Value = "000990.KS"
$ws.Range("D4").Value = 65100
$ws.Range("E4").Value = 41.9
$ws.Range("F4").Value = 7.07
$ws.Range("G4").Value = 40
$ws.Range("H4").Value = 43
$ws.Range("I4").Value = 53
$ws.Range("J4").Value = 60
$ws.Range("K4").Value = 59
$ws.Range("N4").Value = 85.92500513438651

# Row 5 - now SK hynix (000660.KS)
$ws.Range("B5").Value = "SK hynix"
$ws.Range("C5").Value = "000660.KS"
$ws.Range("D5").Value = 544000
$ws.Range("E5").Value = 35.8
$ws.Range("F5").Value = 4.69
$ws.Range("G5").Value = 20
$ws.Range("H5").Value = 53
$ws.Range("I5").Value = 66
$ws.Range("J5").Value = 73
$ws.Range("K5").Value = 58.2
$ws.Range("N5").Value = 85.92500513438651

# Row 6 - now 058470.KS ticker group
$ws.Range("B6").Value = "058470.KS,0P0000ASU1,98886"
$ws.Range("C6").Value = "058470.KS"
$ws.Range("D6").Value = 68300
$ws.Range("E6").Value = 71.40000000000001
$ws.Range("F6").Value = 25.55
$ws.Range("G6").Value = 40
$ws.Range("H6").Value = 40
$ws.Range("I6").Value = 50
$ws.Range("J6").Value = 66
$ws.Range("K6").Value = 57.8
$ws.Range("N6").Value = 85.92500513438651

# Row 7 - 403870.KS ticker group (values unchanged except N)
$ws.Range("N7").Value = 85.92500513438651
